$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-21 (columns A, B, C)
$data = @(
    @("2006_2007", "market",    0.1058),
    @("2006_2007", "rate",      0.0143),
    @("2006_2007", "credit",    0.1327),
    @("2006_2007", "interest", -0.0263),
    @("2006_2007", "inflation",-0.2366),
    @("2008_2009", "credit",   -0.5286999999999999),
    @("2008_2009", "market",    0.2252),
    @("2008_2009", "rate",     -0.6975),
    @("2008_2009", "inflation",-0.1311),
    @("2008_2009", "interest", -0.4211),
    @("2010_2019", "rate",      0.5772),
    @("2010_2019", "credit",   -0.4083),
    @("2010_2019", "market",    0.865),
    @("2010_2019", "interest",  0.4966),
    @("2010_2019", "inflation", 0.7827),
    @("2020_2023", "market",    0.643),
    @("2020_2023", "inflation", 0.5832000000000001),
    @("2020_2023", "rate",      0.7486),
    @("2020_2023", "credit",    0.4103),
    @("2020_2023", "interest",  0.8001)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
